$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    # Force the cell to be stored as text, even if the string looks like a number,
    # then restore the default (unstyled) appearance so no stray formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "41.561.67"
Set-TextValue 2 5 "  +0.22%  "
Set-TextValue 3 4 "2.458.61"
Set-TextValue 3 5 "  +0.87%  "
Set-TextValue 4 5 "  -1.09%  "
Set-TextValue 5 4 "314.44"
Set-TextValue 5 5 "  +1.60%  "
Set-TextValue 6 4 "91.32"
Set-TextValue 6 5 "  +1.70%  "
Set-TextValue 7 5 "  +3.33%  "
Set-TextValue 8 5 "  -1.01%  "
Set-TextValue 9 4 "0.509"
Set-TextValue 9 5 "  +6.02%  "
Set-TextValue 10 4 "32.49"
Set-TextValue 10 5 "  +2.86%  "
Set-TextValue 11 4 "0.0797"
Set-TextValue 11 5 "  +4.14%  "
Set-TextValue 12 5 "  +1.45%  "
Set-TextValue 13 4 "2.838.13"
Set-TextValue 13 5 "  +1.03%  "
Set-TextValue 14 5 "  +2.42%  "
Set-TextValue 15 4 "15.76"
Set-TextValue 15 5 "  +5.65%  "
Set-TextValue 16 4 "2.447.10"
Set-TextValue 16 5 "  -0.58%  "
Set-TextValue 17 5 "  +2.64%  "
Set-TextValue 18 4 "41.556.65"
Set-TextValue 18 5 "  +1.21%  "
Set-TextValue 19 5 "  +5.50%  "
Set-TextValue 20 4 "0.0₃0936"
Set-TextValue 20 5 "  +4.26%  "
Set-TextValue 21 4 "70.87"
Set-TextValue 21 5 "  +3.07%  "
Set-TextValue 22 4 "11.32"
Set-TextValue 22 5 "  +6.21%  "
Set-TextValue 23 4 "236.91"
Set-TextValue 23 5 "  +2.98%  "
Set-TextValue 24 4 "2.71"
Set-TextValue 24 5 "  +2.08%  "
Set-TextValue 25 5 "  -0.11%  "
Set-TextValue 26 5 "  +3.47%  "
Set-TextValue 27 4 "24.27"
Set-TextValue 27 5 "  +3.60%  "
Set-TextValue 28 5 "  +2.66%  "
Set-TextValue 29 5 "  +2.25%  "
Set-TextValue 30 4 "34.97"
Set-TextValue 30 5 "  +1.03%  "
Set-TextValue 31 4 "155.68"
Set-TextValue 31 5 "  +3.03%  "
Set-TextValue 32 5 "  +3.91%  "
Set-TextValue 33 4 "2.57"
Set-TextValue 33 5 "  +2.12%  "
Set-TextValue 34 4 "0.0758"
Set-TextValue 34 5 "  +2.53%  "
Set-TextValue 35 4 "17.48"
Set-TextValue 35 5 "  -0.27%  "
Set-TextValue 36 5 "  -2.80%  "
Set-TextValue 37 4 "2.88"
Set-TextValue 37 5 "  +0.04%  "
Set-TextValue 38 5 "  +3.37%  "
Set-TextValue 39 5 "  +4.20%  "
Set-TextValue 40 5 "  +0.13%  "
Set-TextValue 41 4 "3.93"
Set-TextValue 41 5 "  -1.01%  "
Set-TextValue 42 5 "  -1.39%  "
Set-TextValue 43 4 "1.963.78"
Set-TextValue 43 5 "  +2.69%  "
Set-TextValue 44 5 "  +3.09%  "
Set-TextValue 45 4 "18.47"
Set-TextValue 45 5 "  -2.57%  "
Set-TextValue 46 5 "  +1.98%  "
Set-TextValue 47 4 "8.96"
Set-TextValue 47 5 "  +5.35%  "
Set-TextValue 48 4 "2.697.08"
Set-TextValue 48 5 "  +0.98%  "
Set-TextValue 49 4 "96.32"
Set-TextValue 49 5 "  +3.54%  "
Set-TextValue 50 4 "66.36"
Set-TextValue 50 5 "  +2.18%  "
Set-TextValue 51 2 "Algorand"
Set-TextValue 51 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 51 4 "0.171"
Set-TextValue 51 5 "  +0.02%  "
